# Daily attendance processing - 2025-12-18 16:36:44
#
# Normalizes the "Recorded By" (column G) attendance-audit strings so that
# "System" always appears first among the comma-separated recorder list.
# This is a straight literal re-ordering of a handful of recurring value
# combinations; any value that doesn't match one of these combinations is
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "system, System, backup@backdoor.com" = "System, backup@backdoor.com, system"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "admin@admin.com, System"             = "System, admin@admin.com"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$lastRow = $ws.UsedRange.Rows.Count
$col = 7  # Column G = "Recorded By"

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $current = $cell.Value2
    if ($replacements.ContainsKey($current)) {
        $cell.Value = $replacements[$current]
    }
}
